# Remove the "The log-likelihood is" sentence paragraph together with the
# following centered math paragraph that renders the log-likelihood formula.
# (Commit message: "Updated rendering and script")

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "The log-likelihood is") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The formula paragraph immediately follows the sentence paragraph.
    $formulaPara = $target.Next()

    $startPos = $target.Range.Start
    $endPos = $formulaPara.Range.End

    $killRange = $d.Range($startPos, $endPos)
    $killRange.Delete()
}
